# "se subio las nuevas clases" - add two new class-day columns (K, L) of
# attendance marks ("p") to the existing control_asistencia sheet, mirroring
# the existing columns (E..J) that already carry the same "p" marker.
# Rows 6 and 9 already have an excused/justified mark in another column for
# that day, so only the very last new class day (column L) is filled in for
# those two rows - matching the source data exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$lastRow = 22
$firstRow = 3

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $ws.Range("L$row").Value = "p"
    if ($row -ne 6 -and $row -ne 9) {
        $ws.Range("K$row").Value = "p"
    }
}

# Update the saved cursor position / selection to reflect the last cell
# touched while filling in the new classes.
$ws.Activate()
$ws.Range("L20").Select()

# Best-effort cosmetic touch-ups matching the rest of the author's resave
# (sheet-tab slider ratio, duplicated "Excel Built-in" style prefix from the
# originating app's style-sheet re-export). These are view/metadata-only and
# have no effect on any cell content.
$excel.ActiveWindow.TabRatio = 0.504
$builtinStyle = $wb.Styles.Item("Excel Built-in Normal 2")
$builtinStyle.Name = "Excel Built-in Excel Built-in Excel Built-in Normal 2"
